$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the existing data row above (row 36) and insert the copy at row 37,
# which shifts the blank/summary rows down by one and clones styles +
# (adjusted) shared formulas exactly like a manual "insert copied row" in Excel.
$ws.Rows.Item(36).Copy()
$ws.Rows.Item(37).Insert()

# Overwrite the copied values with the new working-hours entry.
$ws.Range("A37").Value = 2014
$ws.Range("B37").Value = 3
$ws.Range("C37").Value = 2
$ws.Range("D37").Value = 0.6875
$ws.Range("E37").Value = 0.75

# Update the selection to F37 as indicated by the diff.
$ws.Range("F37").Select()
